# edit.ps1 -- apply the "Brugergraenseflade" diff via Word COM-interop
# Strategy: locate each affected paragraph via Find (stable text anchors),
# expand the found range to the whole paragraph (wdParagraph = 4), and
# replace its OOXML with a precisely-built replacement using Range.InsertXML.
# Paragraph deletions are done with Range.Delete(). Operations run from the
# bottom of the document upward so earlier, not-yet-touched text offsets
# used by later Find calls stay valid.

function Get-ParaRange($doc, $searchText) {
    $r = $doc.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $whole = $r.Duplicate
    $whole.Expand(4) | Out-Null
    return $whole
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Figur 2 - ..." caption: convert <w:fldSimple> to begin/instrText/
#    separate/result/end <w:fldChar> runs.
# ------------------------------------------------------------------
$p = Get-ParaRange $d "Navigationsdiagram over Salg"
$p.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008C6B2F" w:rsidRDefault="001A1385" w:rsidP="001A1385"><w:pPr><w:pStyle w:val="Billedtekst"/><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">Figur </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> SEQ Figur \* ARABIC </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve"> - Navigationsdiagram over Salg</w:t></w:r></w:p>')

# ------------------------------------------------------------------
# 2) "Lyn Salg:" list paragraph gains <w:lastRenderedPageBreak/> on its
#    first run.
# ------------------------------------------------------------------
$p = Get-ParaRange $d "Lyn Salg:"
$p.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B5348E" w:rsidRDefault="00B5348E" w:rsidP="00B5348E"><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Lyn Salg:</w:t></w:r><w:r w:rsidR="00310B0A"><w:t xml:space="preserve"> Denn</w:t></w:r><w:r><w:t xml:space="preserve">e bruges ved kontantsalg, altså hvor kunden kommer ned i butikken og vedkommende vælger en eller flere varer fra lageret, som betales og udleveres med det samme i butikken. Der bliver ikke gemt </w:t></w:r><w:r w:rsidR="00310B0A"><w:t>nogen</w:t></w:r><w:r><w:t xml:space="preserve"> informationer om kunden, vedkommende får bare en faktura med som bevis for betaling.</w:t></w:r></w:p>')

# ------------------------------------------------------------------
# 3) "Vi har valgt at lave et navigationsdiagram..." paragraph loses its
#    <w:lastRenderedPageBreak/>.
# ------------------------------------------------------------------
$p = Get-ParaRange $d "Vi har valgt at lave et navigationsdiagram"
$p.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00980230" w:rsidRDefault="00980230"><w:r><w:t>Vi har valgt at lave et navigationsdiagram over salg, der viser navigationen i programmet fra det åbnes til salget er udført.</w:t></w:r><w:r w:rsidR="00B5348E"><w:t xml:space="preserve"> Hele omdrejningspunktet i systemet er salgsmenuen, det er derfor vigtigt at gøre denne del så brugervenlig som overhoved mulig. Vi har valgt at dele salg op i 2 underpunkter </w:t></w:r></w:p>')

# ------------------------------------------------------------------
# 4) Delete the "< Indsaet billeder der viser trin ... >" placeholder
#    paragraph plus the two empty paragraphs right after it.
# ------------------------------------------------------------------
$placeholderPara = Get-ParaRange $d "Indsæt billeder der viser trin"
$delStart = $placeholderPara.Start
$delEnd = $placeholderPara.End
for ($k = 0; $k -lt 2; $k++) {
    $nextPara = $d.Range($delEnd, $delEnd).Paragraphs(1)
    $delEnd = $nextPara.Range.End
}
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete() | Out-Null

# ------------------------------------------------------------------
# 5) "Ydermere har vi ..." paragraph: expand the sentence and move the
#    _GoBack bookmark to sit right before the final sentence.
# ------------------------------------------------------------------
$p = Get-ParaRange $d "Ydermere har vi for at brugeren nemmer"
$p.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00604FAB" w:rsidRDefault="007E4587"><w:r><w:t>Ydermere har vi for at brugeren nemmer</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t xml:space="preserve"> kan bevare overblikket valgt at m</w:t></w:r><w:r><w:t xml:space="preserve">an trinvis udfylder </w:t></w:r><w:r><w:t>ordre</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Så man starter med at udfylde informationer om kunden samt ordren</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> hvorefter ved klik på videre kommer man til hvor man kan tilføje</w:t></w:r><w:r><w:t xml:space="preserve"> vare ting ordren. Dette  er  </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>tydeliggjort i navigationsdiagrammet nedenfor.</w:t></w:r></w:p>')

# ------------------------------------------------------------------
# 6) "Figur 1 - ..." caption: same fldSimple -> fldChar conversion as #1.
# ------------------------------------------------------------------
$p = Get-ParaRange $d "Skitse over Bestilling"
$p.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00980230" w:rsidRDefault="001A1385" w:rsidP="001A1385"><w:pPr><w:pStyle w:val="Billedtekst"/></w:pPr><w:r><w:t xml:space="preserve">Figur </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> SEQ Figur \* ARABIC </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve"> – Skitse over Bestilling</w:t></w:r></w:p>')

Write-Output "done"
